# Regenerate orders with updated distance/size codes.
# Mapping applied to every text cell in the used range:
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31   (S20 / S25 unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -is [string]) {
            $newV = $v -replace 'D64', 'D69'
            $newV = $newV -replace 'D51', 'D55'
            $newV = $newV -replace 'D80', 'D86'
            $newV = $newV -replace 'S30', 'S31'
            if ($newV -ne $v) {
                $cell.Value = $newV
            }
        }
    }
}

Write-Host "Done updating distance/size codes"
